$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1" (sheet1): update the conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 14.36 = 59052.4 pesos`n✅ 59052.4 pesos = 14.32 = 973.71 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas" (sheet2): update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("O10").Value = 4113
$wsTasas.Range("N12").Value = 4123.99
$wsTasas.Range("O12").Value = 68
